$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.584.41'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.024.26'
$ws.Range('E3').Value = '  +2.07%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '379.57'
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.17'
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.588'
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.68'
$ws.Range('E10').Value = '  +0.32%  '
$ws.Range('E11').Value = '  -0.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0859'
$ws.Range('E12').Value = '  +0.66%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.503.83'
$ws.Range('E13').Value = '  +2.33%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.49'
$ws.Range('E14').Value = '  +0.27%  '
$ws.Range('E15').Value = '  -0.38%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.036.73'
$ws.Range('E16').Value = '  +2.35%  '
$ws.Range('E17').Value = '  -3.71%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.60'
$ws.Range('E18').Value = '  -14.63%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '51.592.38'
$ws.Range('E19').Value = '  +0.88%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.09'
$ws.Range('E20').Value = '  +0.41%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('E22').Value = '  +0.17%  '
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '267.23'
$ws.Range('E24').Value = '  -0.75%  '
$ws.Range('E25').Value = '  -6.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.37'
$ws.Range('E26').Value = '  +4.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.51'
$ws.Range('E27').Value = '  +7.46%  '
$ws.Range('E28').Value = '  +3.44%  '
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '26.14'
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.25'
$ws.Range('E32').Value = '  -2.90%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.07'
$ws.Range('E33').Value = '  +0.44%  '
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '34.00'
$ws.Range('E34').Value = '  -0.77%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '50.58'
$ws.Range('E35').Value = '  -1.27%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0450'
$ws.Range('E36').Value = '  +3.24%  '
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.30'
$ws.Range('E38').Value = '  +1.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.299'
$ws.Range('E39').Value = '  +14.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.02'
$ws.Range('E40').Value = '  +1.75%  '
$ws.Range('E41').Value = '  +0.88%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '126.97'
$ws.Range('E42').Value = '  +2.43%  '
$ws.Range('E43').Value = '  -1.20%  '
$ws.Range('E44').Value = '  +1.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.74'
$ws.Range('E45').Value = '  +4.79%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '21.54'
$ws.Range('E46').Value = '  -1.16%  '
$ws.Range('E47').Value = '  +3.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.07'
$ws.Range('E48').Value = '  +1.73%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.018.56'
$ws.Range('E49').Value = '  -3.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.323.40'
$ws.Range('E50').Value = '  +2.14%  '
$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.513'
$ws.Range('E51').Value = '  +4.46%  '
